$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation row needs to be inserted between the existing
# "2026/01/09" rows (610, 611) and the "2026/12/29" block, shifting
# everything from the old row 612 onward down by one row.
$ws.Rows.Item(612).Insert()

# Fill in the newly inserted row. The date column must stay plain text
# (matching the rest of the sheet, which stores dates as literal
# strings rather than date serials), so force a text number format
# before assigning the value and then clear the formatting again so the
# cell keeps the sheet's default (unstyled) look.
$ws.Cells.Item(612, 1).NumberFormat = "@"
$ws.Cells.Item(612, 1).Value = "2026/01/09"
$ws.Cells.Item(612, 1).ClearFormats()

$ws.Cells.Item(612, 2).Value = "金"
$ws.Cells.Item(612, 3).Value = 13
$ws.Cells.Item(612, 4).Value = 180
